$d = $word.ActiveDocument

# 1) UEFA "After 1994" h2h==0 section: the goal-difference-threshold rule was
#    wrong at -1; it should also allow -2 (not just -1), per the commit
#    message ("gls diff can be -2 ... not just -1"). This is the 5th/last
#    occurrence of this sentence in the doc, inside the "pts diff == -2"
#    drawing branch (After 1994). Target it via the specific paragraph so
#    the four earlier, unrelated occurrences of the same sentence are left
#    untouched.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "if gls diff == -1 and gls scored lagging*") {
        $targetPara = $p
    }
}
if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute("gls diff == -1 and gls scored lagging", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "gls diff == -2 and gls scored lagging", 2)
}

# 2) Drop the trailing "Group A 2024" scratch section (three blank
#    paragraphs plus the heading/bullets) that was removed from the end of
#    the document.
$lastGoodIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*leading is winning (last game points == 3), lagging only one points behind*") {
        $lastGoodIndex = $i
    }
}
if ($lastGoodIndex -gt 0 -and $lastGoodIndex -lt $d.Paragraphs.Count) {
    $startPara = $d.Paragraphs.Item($lastGoodIndex + 1)
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $delRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
    $delRange.Delete()
}
